# Fruta / hortaliza, semanal
# Insert two new weekly price rows (date 44628) for Albahaca - Mercado Mayorista
# Lo Valledor de Santiago, right above the former row 413, pushing the
# existing rows 413-419 down to 415-421.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 413 (old rows 413-419 shift to 415-421)
$ws.Rows.Item(413).EntireRow.Insert()
$ws.Rows.Item(413).EntireRow.Insert()

# --- New row 413 ---
$ws.Cells.Item(413,1).Value2 = 6
$ws.Cells.Item(413,2).Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(413,3).Value2 = "Metropolitana"
$ws.Cells.Item(413,4).Value2 = 44628
$ws.Cells.Item(413,5).Value2 = 13
$ws.Cells.Item(413,6).Value2 = 100112052
$ws.Cells.Item(413,7).Value2 = "Albahaca"
$ws.Cells.Item(413,8).Value2 = "Sin especificar"
$ws.Cells.Item(413,9).Value2 = "Primera"
$ws.Cells.Item(413,10).Value2 = 3200
$ws.Cells.Item(413,11).Value2 = 3000
$ws.Cells.Item(413,12).Value2 = 3500
$ws.Cells.Item(413,13).Value2 = 3055
$ws.Cells.Item(413,14).Value2 = "`$/docena de matas"
$ws.Cells.Item(413,15).Value2 = "Región Metropolitana"
$ws.Cells.Item(413,16).Value2 = 509
$ws.Cells.Item(413,17).Value2 = 6
$ws.Cells.Item(413,18).Value2 = "Hortaliza"

# --- New row 414 ---
$ws.Cells.Item(414,1).Value2 = 6
$ws.Cells.Item(414,2).Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(414,3).Value2 = "Metropolitana"
$ws.Cells.Item(414,4).Value2 = 44628
$ws.Cells.Item(414,5).Value2 = 13
$ws.Cells.Item(414,6).Value2 = 100112052
$ws.Cells.Item(414,7).Value2 = "Albahaca"
$ws.Cells.Item(414,8).Value2 = "Sin especificar"
$ws.Cells.Item(414,9).Value2 = "Segunda"
$ws.Cells.Item(414,10).Value2 = 1800
$ws.Cells.Item(414,11).Value2 = 2000
$ws.Cells.Item(414,12).Value2 = 2500
$ws.Cells.Item(414,13).Value2 = 2056
$ws.Cells.Item(414,14).Value2 = "`$/docena de matas"
$ws.Cells.Item(414,15).Value2 = "Región Metropolitana"
$ws.Cells.Item(414,16).Value2 = 343
$ws.Cells.Item(414,17).Value2 = 6
$ws.Cells.Item(414,18).Value2 = "Hortaliza"
